$wb = $excel.ActiveWorkbook

# --- Sheet "Q1_20_21": remove the Mars, A13 and Columbia rows, keeping SoT and F9 ---
$ws1 = $wb.Worksheets.Item("Q1_20_21")
$ws1.Rows.Item(7).Delete()   # Columbia (was row 7)
$ws1.Rows.Item(5).Delete()   # A13 (was row 5)
$ws1.Rows.Item(3).Delete()   # Mars (was row 3)

# --- Sheet "Q4_19_20": remove the Mars, A11 and F9 rows, keeping SoT, A13 and Columbia ---
$ws2 = $wb.Worksheets.Item("Q4_19_20")
$ws2.Rows.Item(7).Delete()   # F9 (was row 7)
$ws2.Rows.Item(5).Delete()   # A11 (was row 5)
$ws2.Rows.Item(3).Delete()   # Mars (was row 3)
